$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 734741.9399999999
$ws.Range("J17").Value = 799528
$ws.Range("L17").Value = 2398584
$ws.Range("N17").Value = -2398920
$ws.Range("H64").Value = 2201466
$ws.Range("I64").Value = 3348081.5
$ws.Range("K64").Value = 3348081.5
$ws.Range("M64").Value = -3347833.5
$ws.Range("H67").Value = 2201466
$ws.Range("I67").Value = 3348081.5
$ws.Range("K67").Value = 3348081.5
$ws.Range("M67").Value = -3347223.5
$ws.Range("H74").Value = 3804.48
$ws.Range("I74").Value = 3242.1177
$ws.Range("J74").Value = 4999.5
$ws.Range("K74").Value = 3242.1177
$ws.Range("L74").Value = 4999.5
$ws.Range("M74").Value = -2306.1177
$ws.Range("N74").Value = -6871.5
$ws.Range("H77").Value = 3804.48
$ws.Range("I77").Value = 3242.1177
$ws.Range("J77").Value = 4999.5
$ws.Range("K77").Value = 16210.5885
$ws.Range("L77").Value = 24997.5
$ws.Range("M77").Value = -11530.5885
$ws.Range("N77").Value = -34357.5
$ws.Range("H107").Value = 8589.076999999999
$ws.Range("I107").Value = 9296
$ws.Range("K107").Value = 9296
$ws.Range("M107").Value = -7376
$ws.Range("H132").Value = 4905845
$ws.Range("I132").Value = 3075.077
$ws.Range("J132").Value = 20839846
$ws.Range("K132").Value = 9225.231
$ws.Range("L132").Value = 62519538
$ws.Range("M132").Value = -6695.231
$ws.Range("N132").Value = -62524598
$ws.Range("H138").Value = 4568401
$ws.Range("I138").Value = 8334498.5
$ws.Range("J138").Value = 3434.0303
$ws.Range("K138").Value = 25003495.5
$ws.Range("L138").Value = 10302.0909
$ws.Range("M138").Value = -24998355.5
$ws.Range("N138").Value = -20582.0909

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3494.23
$ws.Range("I32").Value = 3158.398
$ws.Range("J32").Value = 19950
$ws.Range("K32").Value = 3158.398
$ws.Range("L32").Value = 19950
$ws.Range("M32").Value = -2871.398
$ws.Range("N32").Value = -20524
$ws.Range("H61").Value = 2044.4667
$ws.Range("I61").Value = 2119.6904
$ws.Range("J61").Value = 991.3333
$ws.Range("K61").Value = 2119.6904
$ws.Range("L61").Value = 991.3333
$ws.Range("M61").Value = -1907.6904
$ws.Range("N61").Value = -1415.3333
$ws.Range("H74").Value = 1794.641
$ws.Range("I74").Value = 1130.2693
$ws.Range("K74").Value = 1130.2693
$ws.Range("M74").Value = -256.2692999999999
$ws.Range("H77").Value = 1794.641
$ws.Range("I77").Value = 1130.2693
$ws.Range("K77").Value = 5651.3465
$ws.Range("M77").Value = -1283.3465
$ws.Range("H132").Value = 3624593.8
$ws.Range("I132").Value = 4718232.5
$ws.Range("K132").Value = 14154697.5
$ws.Range("M132").Value = -14152167.5
$ws.Range("H136").Value = 2044.4667
$ws.Range("I136").Value = 2119.6904
$ws.Range("J136").Value = 991.3333
$ws.Range("K136").Value = 6359.0712
$ws.Range("L136").Value = 2973.9999
$ws.Range("M136").Value = -3809.0712
$ws.Range("N136").Value = -8073.9999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1881.6923
$ws.Range("I31").Value = 1516.0857
$ws.Range("J31").Value = 2179.279
$ws.Range("K31").Value = 1516.0857
$ws.Range("L31").Value = 2179.279
$ws.Range("M31").Value = -1221.0857
$ws.Range("N31").Value = -2769.279
$ws.Range("H34").Value = 1881.6923
$ws.Range("I34").Value = 1516.0857
$ws.Range("J34").Value = 2179.279
$ws.Range("K34").Value = 1516.0857
$ws.Range("L34").Value = 2179.279
$ws.Range("M34").Value = -1314.0857
$ws.Range("N34").Value = -2583.279
$ws.Range("H86").Value = 4307.1177
$ws.Range("I86").Value = 4120.1665
$ws.Range("J86").Value = 4409.091
$ws.Range("K86").Value = 4120.1665
$ws.Range("L86").Value = 4409.091
$ws.Range("M86").Value = -2997.1665
$ws.Range("N86").Value = -6655.091
$ws.Range("H89").Value = 4307.1177
$ws.Range("I89").Value = 4120.1665
$ws.Range("J89").Value = 4409.091
$ws.Range("K89").Value = 20600.8325
$ws.Range("L89").Value = 22045.455
$ws.Range("M89").Value = -14984.8325
$ws.Range("N89").Value = -33277.455
$ws.Range("H99").Value = 1999.6666
$ws.Range("I99").Value = 1999.6666
$ws.Range("K99").Value = 1999.6666
$ws.Range("M99").Value = -501.6666
$ws.Range("H122").Value = 1294.5714
$ws.Range("I122").Value = 1365.1111
$ws.Range("J122").Value = 871.3333
$ws.Range("K122").Value = 4095.3333
$ws.Range("L122").Value = 2613.9999
$ws.Range("M122").Value = -1645.3333
$ws.Range("N122").Value = -7513.9999
$ws.Range("H126").Value = 1999.6666
$ws.Range("I126").Value = 1999.6666
$ws.Range("K126").Value = 5998.9998
$ws.Range("M126").Value = -3528.9998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1107.6666
$ws.Range("I34").Value = 161.66667
$ws.Range("K34").Value = 485.00001
$ws.Range("M34").Value = -401.00001
$ws.Range("H39").Value = 501.1905
$ws.Range("J39").Value = 501.2195
$ws.Range("L39").Value = 1503.6585
$ws.Range("N39").Value = -2091.6585
$ws.Range("H55").Value = 448.5
$ws.Range("I55").Value = 152.33333
$ws.Range("J55").Value = 547.2222
$ws.Range("K55").Value = 456.99999
$ws.Range("L55").Value = 1641.6666
$ws.Range("M55").Value = -279.99999
$ws.Range("N55").Value = -1995.6666
$ws.Range("H129").Value = 2536.6296
$ws.Range("I129").Value = 613.5714
$ws.Range("J129").Value = 4607.615
$ws.Range("K129").Value = 1840.7142
$ws.Range("L129").Value = 13822.845
$ws.Range("M129").Value = 3159.2858
$ws.Range("N129").Value = -23822.845

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11194.192
$ws.Range("I70").Value = 21150.908
$ws.Range("K70").Value = 21150.908
$ws.Range("M70").Value = -20880.908
$ws.Range("H73").Value = 11194.192
$ws.Range("I73").Value = 21150.908
$ws.Range("K73").Value = 21150.908
$ws.Range("M73").Value = -20214.908
$ws.Range("H132").Value = 4113.0713
$ws.Range("I132").Value = 4273.75
$ws.Range("J132").Value = 3149
$ws.Range("K132").Value = 12821.25
$ws.Range("L132").Value = 9447
$ws.Range("M132").Value = -10291.25
$ws.Range("N132").Value = -14507

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4985.971
$ws.Range("I7").Value = 5376.8125
$ws.Range("J7").Value = 4656.8423
$ws.Range("K7").Value = 5376.8125
$ws.Range("L7").Value = 4656.8423
$ws.Range("M7").Value = -5264.8125
$ws.Range("N7").Value = -4880.8423
$ws.Range("H40").Value = 3752.9211
$ws.Range("I40").Value = 4517.647
$ws.Range("J40").Value = 3133.8572
$ws.Range("K40").Value = 4517.647
$ws.Range("L40").Value = 3133.8572
$ws.Range("M40").Value = -4381.647
$ws.Range("N40").Value = -3405.8572
$ws.Range("H61").Value = 1962.9412
$ws.Range("I61").Value = 1812.2858
$ws.Range("J61").Value = 2666
$ws.Range("K61").Value = 1812.2858
$ws.Range("L61").Value = 2666
$ws.Range("M61").Value = -1610.2858
$ws.Range("N61").Value = -3070
$ws.Range("H68").Value = 2365.1538
$ws.Range("I68").Value = 2337.1428
$ws.Range("J68").Value = 2397.8333
$ws.Range("K68").Value = 2337.1428
$ws.Range("L68").Value = 2397.8333
$ws.Range("M68").Value = -1588.1428
$ws.Range("N68").Value = -3895.8333
$ws.Range("H71").Value = 2365.1538
$ws.Range("I71").Value = 2337.1428
$ws.Range("J71").Value = 2397.8333
$ws.Range("K71").Value = 11685.714
$ws.Range("L71").Value = 11989.1665
$ws.Range("M71").Value = -7941.714
$ws.Range("N71").Value = -19477.1665
$ws.Range("H82").Value = 1781.4166
$ws.Range("I82").Value = 1663.7778
$ws.Range("J82").Value = 2134.3333
$ws.Range("K82").Value = 1663.7778
$ws.Range("L82").Value = 2134.3333
$ws.Range("M82").Value = -1302.7778
$ws.Range("N82").Value = -2856.3333
$ws.Range("H85").Value = 1781.4166
$ws.Range("I85").Value = 1663.7778
$ws.Range("J85").Value = 2134.3333
$ws.Range("K85").Value = 1663.7778
$ws.Range("L85").Value = 2134.3333
$ws.Range("M85").Value = -415.7778000000001
$ws.Range("N85").Value = -4630.3333
$ws.Range("H100").Value = 2666.6667
$ws.Range("I100").Value = 2500
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2500
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1959
$ws.Range("N100").Value = -4082
$ws.Range("H113").Value = 1962.9412
$ws.Range("I113").Value = 1812.2858
$ws.Range("J113").Value = 2666
$ws.Range("K113").Value = 1812.2858
$ws.Range("L113").Value = 2666
$ws.Range("M113").Value = 357.7141999999999
$ws.Range("N113").Value = -7006
$ws.Range("H126").Value = 4985.971
$ws.Range("I126").Value = 5376.8125
$ws.Range("J126").Value = 4656.8423
$ws.Range("K126").Value = 16130.4375
$ws.Range("L126").Value = 13970.5269
$ws.Range("M126").Value = -13660.4375
$ws.Range("N126").Value = -18910.5269
$ws.Range("H132").Value = 9265017
$ws.Range("I132").Value = 3444.9512
$ws.Range("J132").Value = 38474590
$ws.Range("K132").Value = 10334.8536
$ws.Range("L132").Value = 115423770
$ws.Range("M132").Value = -7804.8536
$ws.Range("N132").Value = -115428830

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1020.73334
$ws.Range("I107").Value = 1356.7778
$ws.Range("J107").Value = 516.6667
$ws.Range("K107").Value = 4070.3334
$ws.Range("L107").Value = 1550.0001
$ws.Range("M107").Value = -2150.3334
$ws.Range("N107").Value = -5390.0001
$ws.Range("H122").Value = 1990.875
$ws.Range("I122").Value = 1961.4615
$ws.Range("J122").Value = 2118.3333
$ws.Range("K122").Value = 5884.3845
$ws.Range("L122").Value = 6354.999899999999
$ws.Range("M122").Value = -3434.3845
$ws.Range("N122").Value = -11254.9999
$ws.Range("H126").Value = 1505.0385
$ws.Range("I126").Value = 1570.9131
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 4712.7393
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -2242.7393
$ws.Range("N126").Value = -7940
